$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correction of English text: the column header in C1 ("mark") is renamed
# to the more standard English term "grade".
$ws.Range("C1").Value = "grade"

# The active cell / selection moves from A10 to A12.
$ws.Range("A12").Select()
